# Auto commit 28-05-2025 11:43:03.83
# Rename the borewell/tubewell construction sheets to their full names and
# update each sheet's remembered selection / scroll position, then restore
# the originally-intended active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$wb.Worksheets.Item("110 mm BWC").Name = "110 mm Borewell Construction"
$wb.Worksheets.Item("150 mm BWC").Name = "150 mm Borewell Construction"
$wb.Worksheets.Item("150 mm TWC").Name = "150 mm Tubewell Construction"
$wb.Worksheets.Item("200 mm TWC").Name = "200 mm Tubewell Construction"

# --- 2. Update each sheet's remembered selection / scroll position ----

# "150 mm Borewell Construction": scrolled so row 7 is at the top, I21 selected
$ws = $wb.Worksheets.Item("150 mm Borewell Construction")
$ws.Activate()
$ws.Range("I21").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# "150 mm Tubewell Construction": scrolled so row 16 is at the top, J20 selected
$ws = $wb.Worksheets.Item("150 mm Tubewell Construction")
$ws.Activate()
$ws.Range("J20").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1

# "200 mm Tubewell Construction": H24 selected
$ws = $wb.Worksheets.Item("200 mm Tubewell Construction")
$ws.Activate()
$ws.Range("H24").Select()

# "Compressor Pump Installation" keeps its E5 selection (no change needed),
# it just stops being the active tab once another sheet is activated below.

# --- 3. Make "110 mm Borewell Construction" the active sheet/tab ------
# This sheet becomes tabSelected="1" / the workbook's active tab (matching
# activeTab moving to this sheet), with F21 as the selected cell.
$ws = $wb.Worksheets.Item("110 mm Borewell Construction")
$ws.Activate()
$ws.Range("F21").Select()
